$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 4823  # H29: 5217.25 -> 4823
$ws.Cells.Item(29, 10).Value = 4823  # J29: 5217.25 -> 4823
$ws.Cells.Item(29, 12).Value = 14469  # L29: 15651.75 -> 14469
$ws.Cells.Item(29, 14).Value = -15031  # N29: -16213.75 -> -15031
$ws.Cells.Item(51, 8).Value = 7322  # H51: 6769.25 -> 7322
$ws.Cells.Item(51, 9).Value = 6481.3335  # I51: 5972 -> 6481.3335
$ws.Cells.Item(51, 10).Value = 8583  # J51: 7566.5 -> 8583
$ws.Cells.Item(51, 11).Value = 6481.3335  # K51: 5972 -> 6481.3335
$ws.Cells.Item(51, 12).Value = 8583  # L51: 7566.5 -> 8583
$ws.Cells.Item(51, 13).Value = -5997.3335  # M51: -5488 -> -5997.3335
$ws.Cells.Item(51, 14).Value = -9551  # N51: -8534.5 -> -9551
$ws.Cells.Item(74, 8).Value = 11364.909  # H74: 10031.654 -> 11364.909
$ws.Cells.Item(74, 9).Value = 13626.75  # I74: 11838.053 -> 13626.75
$ws.Cells.Item(74, 10).Value = 5333.3335  # J74: 5128.5713 -> 5333.3335
$ws.Cells.Item(74, 11).Value = 13626.75  # K74: 11838.053 -> 13626.75
$ws.Cells.Item(74, 12).Value = 5333.3335  # L74: 5128.5713 -> 5333.3335
$ws.Cells.Item(74, 13).Value = -12690.75  # M74: -10902.053 -> -12690.75
$ws.Cells.Item(74, 14).Value = -7205.3335  # N74: -7000.5713 -> -7205.3335
$ws.Cells.Item(77, 8).Value = 11364.909  # H77: 10031.654 -> 11364.909
$ws.Cells.Item(77, 9).Value = 13626.75  # I77: 11838.053 -> 13626.75
$ws.Cells.Item(77, 10).Value = 5333.3335  # J77: 5128.5713 -> 5333.3335
$ws.Cells.Item(77, 11).Value = 68133.75  # K77: 59190.265 -> 68133.75
$ws.Cells.Item(77, 12).Value = 26666.6675  # L77: 25642.8565 -> 26666.6675
$ws.Cells.Item(77, 13).Value = -63453.75  # M77: -54510.265 -> -63453.75
$ws.Cells.Item(77, 14).Value = -36026.6675  # N77: -35002.85649999999 -> -36026.6675
$ws.Cells.Item(113, 8).Value = 33336164  # H113: 33336184 -> 33336164
$ws.Cells.Item(113, 9).Value = 83335630  # I113: 83335690 -> 83335630
$ws.Cells.Item(113, 11).Value = 83335630  # K113: 83335690 -> 83335630
$ws.Cells.Item(113, 13).Value = -83332376  # M113: -83332436 -> -83332376
$ws.Cells.Item(115, 8).Value = 1122.875  # H115: 1123.125 -> 1122.875
$ws.Cells.Item(115, 9).Value = 283.2857  # I115: 283.57144 -> 283.2857
$ws.Cells.Item(115, 11).Value = 849.8571000000001  # K115: 850.71432 -> 849.8571000000001
$ws.Cells.Item(115, 13).Value = 717.1428999999999  # M115: 716.28568 -> 717.1428999999999
$ws.Cells.Item(132, 8).Value = 5147.919  # H132: 5068.8687 -> 5147.919
$ws.Cells.Item(132, 9).Value = 5219.9707  # I132: 5132.086 -> 5219.9707
$ws.Cells.Item(132, 11).Value = 15659.9121  # K132: 15396.258 -> 15659.9121
$ws.Cells.Item(132, 13).Value = -13129.9121  # M132: -12866.258 -> -13129.9121
$ws.Cells.Item(137, 8).Value = 2006918.8  # H137: 2280389.5 -> 2006918.8
$ws.Cells.Item(137, 9).Value = 10001139  # I137: 25000650 -> 10001139
$ws.Cells.Item(137, 11).Value = 30003417  # K137: 75001950 -> 30003417
$ws.Cells.Item(137, 13).Value = -30000867  # M137: -74999400 -> -30000867
$ws.Cells.Item(138, 8).Value = 3829.484  # H138: 3811.0896 -> 3829.484
$ws.Cells.Item(138, 10).Value = 3494.898  # J138: 3503.0557 -> 3494.898
$ws.Cells.Item(138, 12).Value = 10484.694  # L138: 10509.1671 -> 10484.694
$ws.Cells.Item(138, 14).Value = -20764.694  # N138: -20789.1671 -> -20764.694
$ws.Cells.Item(141, 8).Value = 5982.294  # H141: 6330.067 -> 5982.294
$ws.Cells.Item(141, 9).Value = 2909.0833  # I141: 2816.1 -> 2909.0833
$ws.Cells.Item(141, 11).Value = 8727.249899999999  # K141: 8448.299999999999 -> 8727.249899999999
$ws.Cells.Item(141, 13).Value = -3547.249899999999  # M141: -3268.299999999999 -> -3547.249899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1021.7083  # H2: 1087.7273 -> 1021.7083
$ws.Cells.Item(2, 9).Value = 825.1177  # I2: 895.73334 -> 825.1177
$ws.Cells.Item(2, 11).Value = 825.1177  # K2: 895.73334 -> 825.1177
$ws.Cells.Item(2, 13).Value = -712.1177  # M2: -782.73334 -> -712.1177
$ws.Cells.Item(32, 8).Value = 6435.8  # H32: 6238.3076 -> 6435.8
$ws.Cells.Item(32, 9).Value = 5756.952  # I32: 5554.409 -> 5756.952
$ws.Cells.Item(32, 11).Value = 5756.952  # K32: 5554.409 -> 5756.952
$ws.Cells.Item(32, 13).Value = -5469.952  # M32: -5267.409 -> -5469.952
$ws.Cells.Item(44, 8).Value = 62548.5  # H44: 62599 -> 62548.5
$ws.Cells.Item(44, 10).Value = 62548.5  # J44: 62599 -> 62548.5
$ws.Cells.Item(44, 12).Value = 62548.5  # L44: 62599 -> 62548.5
$ws.Cells.Item(44, 14).Value = -63524.5  # N44: -63575 -> -63524.5
$ws.Cells.Item(45, 8).Value = 214210  # H45: 210960 -> 214210
$ws.Cells.Item(45, 9).Value = 420420  # I45: 210960 -> 420420
$ws.Cells.Item(45, 10).Value = 8000  # J45: 0 -> 8000
$ws.Cells.Item(45, 11).Value = 420420  # K45: 210960 -> 420420
$ws.Cells.Item(45, 12).Value = 8000  # L45: 0 -> 8000
$ws.Cells.Item(45, 13).Value = -420043  # M45: -210583 -> -420043
$ws.Cells.Item(45, 14).Value = -8754  # N45: None -> -8754
$ws.Cells.Item(61, 8).Value = 3372.6667  # H61: 3706.7917 -> 3372.6667
$ws.Cells.Item(61, 9).Value = 2416.6  # I61: 2845.8333 -> 2416.6
$ws.Cells.Item(61, 11).Value = 2416.6  # K61: 2845.8333 -> 2416.6
$ws.Cells.Item(61, 13).Value = -2204.6  # M61: -2633.8333 -> -2204.6
$ws.Cells.Item(110, 8).Value = 4238.3335  # H110: 5285 -> 4238.3335
$ws.Cells.Item(110, 9).Value = 1693.8334  # I110: 1917.6 -> 1693.8334
$ws.Cells.Item(110, 10).Value = 9327.333000000001  # J110: 9494.25 -> 9327.333000000001
$ws.Cells.Item(110, 11).Value = 1693.8334  # K110: 1917.6 -> 1693.8334
$ws.Cells.Item(110, 12).Value = 9327.333000000001  # L110: 9494.25 -> 9327.333000000001
$ws.Cells.Item(110, 13).Value = 351.1666  # M110: 127.4000000000001 -> 351.1666
$ws.Cells.Item(110, 14).Value = -13417.333  # N110: -13584.25 -> -13417.333
$ws.Cells.Item(116, 8).Value = 1021.7083  # H116: 1087.7273 -> 1021.7083
$ws.Cells.Item(116, 9).Value = 825.1177  # I116: 895.73334 -> 825.1177
$ws.Cells.Item(116, 11).Value = 825.1177  # K116: 895.73334 -> 825.1177
$ws.Cells.Item(116, 13).Value = 1468.8823  # M116: 1398.26666 -> 1468.8823
$ws.Cells.Item(132, 8).Value = 4833328  # H132: 4904392 -> 4833328
$ws.Cells.Item(132, 9).Value = 2140.9465  # I132: 2161.8728 -> 2140.9465
$ws.Cells.Item(132, 11).Value = 6422.8395  # K132: 6485.6184 -> 6422.8395
$ws.Cells.Item(132, 13).Value = -3892.8395  # M132: -3955.6184 -> -3892.8395
$ws.Cells.Item(136, 8).Value = 3372.6667  # H136: 3706.7917 -> 3372.6667
$ws.Cells.Item(136, 9).Value = 2416.6  # I136: 2845.8333 -> 2416.6
$ws.Cells.Item(136, 11).Value = 7249.799999999999  # K136: 8537.499899999999 -> 7249.799999999999
$ws.Cells.Item(136, 13).Value = -4699.799999999999  # M136: -5987.499899999999 -> -4699.799999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1021.7083  # H3: 1087.7273 -> 1021.7083
$ws.Cells.Item(3, 9).Value = 825.1177  # I3: 895.73334 -> 825.1177
$ws.Cells.Item(3, 11).Value = 825.1177  # K3: 895.73334 -> 825.1177
$ws.Cells.Item(3, 13).Value = -711.1177  # M3: -781.73334 -> -711.1177
$ws.Cells.Item(22, 8).Value = 1000  # H22: 0 -> 1000
$ws.Cells.Item(22, 10).Value = 1000  # J22: 0 -> 1000
$ws.Cells.Item(22, 12).Value = 1000  # L22: 0 -> 1000
$ws.Cells.Item(22, 14).Value = -1346  # N22: None -> -1346

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6164.1816  # H31: 6339.294 -> 6164.1816
$ws.Cells.Item(31, 9).Value = 4520.2856  # I31: 4613.909 -> 4520.2856
$ws.Cells.Item(31, 10).Value = 9041  # J31: 9502.5 -> 9041
$ws.Cells.Item(31, 11).Value = 4520.2856  # K31: 4613.909 -> 4520.2856
$ws.Cells.Item(31, 12).Value = 9041  # L31: 9502.5 -> 9041
$ws.Cells.Item(31, 13).Value = -4225.2856  # M31: -4318.909 -> -4225.2856
$ws.Cells.Item(31, 14).Value = -9631  # N31: -10092.5 -> -9631
$ws.Cells.Item(34, 8).Value = 6164.1816  # H34: 6339.294 -> 6164.1816
$ws.Cells.Item(34, 9).Value = 4520.2856  # I34: 4613.909 -> 4520.2856
$ws.Cells.Item(34, 10).Value = 9041  # J34: 9502.5 -> 9041
$ws.Cells.Item(34, 11).Value = 4520.2856  # K34: 4613.909 -> 4520.2856
$ws.Cells.Item(34, 12).Value = 9041  # L34: 9502.5 -> 9041
$ws.Cells.Item(34, 13).Value = -4318.2856  # M34: -4411.909 -> -4318.2856
$ws.Cells.Item(34, 14).Value = -9445  # N34: -9906.5 -> -9445
$ws.Cells.Item(58, 8).Value = 2686.2856  # H58: 2672.6 -> 2686.2856
$ws.Cells.Item(58, 9).Value = 2167.5557  # I58: 2201 -> 2167.5557
$ws.Cells.Item(58, 10).Value = 3620  # J58: 3380 -> 3620
$ws.Cells.Item(58, 11).Value = 2167.5557  # K58: 2201 -> 2167.5557
$ws.Cells.Item(58, 12).Value = 3620  # L58: 3380 -> 3620
$ws.Cells.Item(58, 13).Value = -1964.5557  # M58: -1998 -> -1964.5557
$ws.Cells.Item(58, 14).Value = -4026  # N58: -3786 -> -4026
$ws.Cells.Item(99, 8).Value = 9553.727999999999  # H99: 9015.833000000001 -> 9553.727999999999
$ws.Cells.Item(99, 10).Value = 3549.5  # J99: 3399.3333 -> 3549.5
$ws.Cells.Item(99, 12).Value = 3549.5  # L99: 3399.3333 -> 3549.5
$ws.Cells.Item(99, 14).Value = -6545.5  # N99: -6395.3333 -> -6545.5
$ws.Cells.Item(126, 8).Value = 9553.727999999999  # H126: 9015.833000000001 -> 9553.727999999999
$ws.Cells.Item(126, 10).Value = 3549.5  # J126: 3399.3333 -> 3549.5
$ws.Cells.Item(126, 12).Value = 10648.5  # L126: 10197.9999 -> 10648.5
$ws.Cells.Item(126, 14).Value = -15588.5  # N126: -15137.9999 -> -15588.5
$ws.Cells.Item(132, 8).Value = 9526979  # H132: 10104326 -> 9526979
$ws.Cells.Item(132, 9).Value = 2678.0476  # I132: 2786.4 -> 2678.0476
$ws.Cells.Item(132, 10).Value = 23813430  # J132: 25645154 -> 23813430
$ws.Cells.Item(132, 11).Value = 8034.1428  # K132: 8359.200000000001 -> 8034.1428
$ws.Cells.Item(132, 12).Value = 71440290  # L132: 76935462 -> 71440290
$ws.Cells.Item(132, 13).Value = -5504.1428  # M132: -5829.200000000001 -> -5504.1428
$ws.Cells.Item(132, 14).Value = -71445350  # N132: -76940522 -> -71445350
$ws.Cells.Item(134, 8).Value = 3895.8667  # H134: 3995.6428 -> 3895.8667
$ws.Cells.Item(134, 9).Value = 3879.8462  # I134: 3994.9167 -> 3879.8462
$ws.Cells.Item(134, 11).Value = 11639.5386  # K134: 11984.7501 -> 11639.5386
$ws.Cells.Item(134, 13).Value = -9104.5386  # M134: -9449.750100000001 -> -9104.5386
$ws.Cells.Item(135, 8).Value = 69750.25  # H135: 75000 -> 69750.25
$ws.Cells.Item(135, 10).Value = 69750.25  # J135: 75000 -> 69750.25
$ws.Cells.Item(135, 12).Value = 69750.25  # L135: 75000 -> 69750.25
$ws.Cells.Item(135, 14).Value = -79890.25  # N135: -85140 -> -79890.25
$ws.Cells.Item(136, 8).Value = 2686.2856  # H136: 2672.6 -> 2686.2856
$ws.Cells.Item(136, 9).Value = 2167.5557  # I136: 2201 -> 2167.5557
$ws.Cells.Item(136, 10).Value = 3620  # J136: 3380 -> 3620
$ws.Cells.Item(136, 11).Value = 6502.6671  # K136: 6603 -> 6502.6671
$ws.Cells.Item(136, 12).Value = 10860  # L136: 10140 -> 10860
$ws.Cells.Item(136, 13).Value = -3952.6671  # M136: -4053 -> -3952.6671
$ws.Cells.Item(136, 14).Value = -15960  # N136: -15240 -> -15960

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 1478.1428  # H2: 1350.7391 -> 1478.1428
$ws.Cells.Item(2, 9).Value = 71.833336  # I2: 63.42857 -> 71.833336
$ws.Cells.Item(2, 11).Value = 431.000016  # K2: 380.57142 -> 431.000016
$ws.Cells.Item(2, 13).Value = -318.000016  # M2: -267.57142 -> -318.000016
$ws.Cells.Item(44, 8).Value = 359.8  # H44: 433.33334 -> 359.8
$ws.Cells.Item(44, 9).Value = 299.66666  # I44: 300 -> 299.66666
$ws.Cells.Item(44, 10).Value = 450  # J44: 540 -> 450
$ws.Cells.Item(44, 11).Value = 898.9999799999999  # K44: 900 -> 898.9999799999999
$ws.Cells.Item(44, 12).Value = 1350  # L44: 1620 -> 1350
$ws.Cells.Item(44, 13).Value = -500.9999799999999  # M44: -502 -> -500.9999799999999
$ws.Cells.Item(44, 14).Value = -2146  # N44: -2416 -> -2146
$ws.Cells.Item(131, 8).Value = 10206.823  # H131: 9763.223 -> 10206.823
$ws.Cells.Item(131, 10).Value = 2273.5  # J131: 2267.7778 -> 2273.5
$ws.Cells.Item(131, 12).Value = 6820.5  # L131: 6803.3334 -> 6820.5
$ws.Cells.Item(131, 14).Value = -16900.5  # N131: -16883.3334 -> -16900.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3153.9092  # H122: 3204.093 -> 3153.9092
$ws.Cells.Item(122, 9).Value = 2435.7942  # I122: 2479.4243 -> 2435.7942
$ws.Cells.Item(122, 11).Value = 7307.382599999999  # K122: 7438.2729 -> 7307.382599999999
$ws.Cells.Item(122, 13).Value = -4857.382599999999  # M122: -4988.2729 -> -4857.382599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3835.2104  # H7: 3748.4 -> 3835.2104
$ws.Cells.Item(7, 9).Value = 3923.9092  # I7: 3771.8333 -> 3923.9092
$ws.Cells.Item(7, 11).Value = 3923.9092  # K7: 3771.8333 -> 3923.9092
$ws.Cells.Item(7, 13).Value = -3811.9092  # M7: -3659.8333 -> -3811.9092
$ws.Cells.Item(41, 8).Value = 0  # H41: 30000 -> 0
$ws.Cells.Item(41, 10).Value = 0  # J41: 30000 -> 0
$ws.Cells.Item(41, 12).Value = 0  # L41: 30000 -> 0
$ws.Cells.Item(41, 14).ClearContents()  # N41: remove (was -30876)
$ws.Cells.Item(46, 8).Value = 2349.1667  # H46: 2199.4285 -> 2349.1667
$ws.Cells.Item(46, 9).Value = 2424.5  # I46: 2199.5 -> 2424.5
$ws.Cells.Item(46, 10).Value = 2198.5  # J46: 2199 -> 2198.5
$ws.Cells.Item(46, 11).Value = 2424.5  # K46: 2199.5 -> 2424.5
$ws.Cells.Item(46, 12).Value = 2198.5  # L46: 2199 -> 2198.5
$ws.Cells.Item(46, 13).Value = -2236.5  # M46: -2011.5 -> -2236.5
$ws.Cells.Item(46, 14).Value = -2574.5  # N46: -2575 -> -2574.5
$ws.Cells.Item(68, 8).Value = 2070.0967  # H68: 2005.4688 -> 2070.0967
$ws.Cells.Item(68, 9).Value = 2056.1785  # I68: 1985.3448 -> 2056.1785
$ws.Cells.Item(68, 11).Value = 2056.1785  # K68: 1985.3448 -> 2056.1785
$ws.Cells.Item(68, 13).Value = -1307.1785  # M68: -1236.3448 -> -1307.1785
$ws.Cells.Item(71, 8).Value = 2070.0967  # H71: 2005.4688 -> 2070.0967
$ws.Cells.Item(71, 9).Value = 2056.1785  # I71: 1985.3448 -> 2056.1785
$ws.Cells.Item(71, 11).Value = 10280.8925  # K71: 9926.724 -> 10280.8925
$ws.Cells.Item(71, 13).Value = -6536.8925  # M71: -6182.724 -> -6536.8925
$ws.Cells.Item(102, 8).Value = 43666  # H102: 44999.332 -> 43666
$ws.Cells.Item(102, 10).Value = 43666  # J102: 44999.332 -> 43666
$ws.Cells.Item(102, 12).Value = 43666  # L102: 44999.332 -> 43666
$ws.Cells.Item(102, 14).Value = -50156  # N102: -51489.332 -> -50156
$ws.Cells.Item(126, 8).Value = 3835.2104  # H126: 3748.4 -> 3835.2104
$ws.Cells.Item(126, 9).Value = 3923.9092  # I126: 3771.8333 -> 3923.9092
$ws.Cells.Item(126, 11).Value = 11771.7276  # K126: 11315.4999 -> 11771.7276
$ws.Cells.Item(126, 13).Value = -9301.7276  # M126: -8845.499899999999 -> -9301.7276
$ws.Cells.Item(132, 8).Value = 3996.6099  # H132: 4061.525 -> 3996.6099
$ws.Cells.Item(132, 10).Value = 5824.625  # J132: 6119.6 -> 5824.625
$ws.Cells.Item(132, 12).Value = 17473.875  # L132: 18358.8 -> 17473.875
$ws.Cells.Item(132, 14).Value = -22533.875  # N132: -23418.8 -> -22533.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(102, 8).Value = 30000  # H102: 32000 -> 30000
$ws.Cells.Item(102, 10).Value = 30000  # J102: 32000 -> 30000
$ws.Cells.Item(102, 12).Value = 30000  # L102: 32000 -> 30000
$ws.Cells.Item(102, 14).Value = -36490  # N102: -38490 -> -36490
$ws.Cells.Item(132, 8).Value = 3997.2273  # H132: 4524.6665 -> 3997.2273
$ws.Cells.Item(132, 9).Value = 4149  # I132: 5827.7144 -> 4149
$ws.Cells.Item(132, 10).Value = 3815.1  # J132: 3695.4546 -> 3815.1
$ws.Cells.Item(132, 11).Value = 12447  # K132: 17483.1432 -> 12447
$ws.Cells.Item(132, 12).Value = 11445.3  # L132: 11086.3638 -> 11445.3
$ws.Cells.Item(132, 13).Value = -9917  # M132: -14953.1432 -> -9917
$ws.Cells.Item(132, 14).Value = -16505.3  # N132: -16146.3638 -> -16505.3
